# Append: 2026-02-03 06:57 JST
# Update the "取得日時" (retrieved datetime) column A for all existing data
# rows on the "ランサーズ" sheet from the previous timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-03 06:57:50"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
